# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.771.32"
$ws.Cells.Item(2, 5).Value = "  +2.23%  "
$ws.Cells.Item(3, 4).Value = "2.111.23"
$ws.Cells.Item(3, 5).Value = "  +8.69%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.29%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "333.88"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +4.19%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.22%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5279"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +3.74%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.4382"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +8.17%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.09029"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +7.32%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "45.88"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +8.65%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "1.178"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +4.89%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "24.94"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.71%  "
$ws.Cells.Item(13, 4).Value = "2.110.14"
$ws.Cells.Item(13, 5).Value = "  +9.38%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.759"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +5.19%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.824"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +6.92%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "97.22"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +4.35%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.12%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001127"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.57%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06672"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +2.45%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "19.10"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.83%  "
$ws.Cells.Item(21, 5).Value = "  +0.24%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "6.358"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +6.17%  "
$ws.Cells.Item(23, 4).Value = "30.853.98"
$ws.Cells.Item(23, 5).Value = "  +2.49%  "
$ws.Cells.Item(24, 5).Value = "  +6.65%  "
$ws.Cells.Item(25, 4).Value = "2.358.74"
$ws.Cells.Item(25, 5).Value = "  +9.79%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.262"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +3.32%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "22.79"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.07%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "2.553"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +9.72%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "162.40"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.28%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "132.87"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.63%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.171"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.97%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.1072"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.35%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "6.227"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +3.54%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "4.052"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.540"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +22.11%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.02601"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.68%  "
$ws.Cells.Item(37, 5).Value = "  +3.59%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.06741"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.92%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "9.520"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +9.20%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "12.71"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +8.16%  "
$ws.Cells.Item(41, 5).Value = "  +5.12%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.6831"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +4.72%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "1.248"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.01%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.6457"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +5.97%  "
$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "14.14"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +5.50%  "
$ws.Cells.Item(46, 2).Value = "Frax"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.50%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.233"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.07%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "3.675"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.27%  "
$ws.Cells.Item(49, 5).Value = "  +5.06%  "
$ws.Cells.Item(50, 5).Value = "  +5.11%  "
$ws.Cells.Item(51, 2).Value = "Quant"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "119.40"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.74%  "
